# Update latest output (run 174)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" (sheet1): rows 3-4, columns E (Cost) and F (Unit Cost) ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E3").Value = 750.4493392500002
$schedule.Range("F3").Value = 28.36165303287983
$schedule.Range("E4").Value = 404.7445155
$schedule.Range("F4").Value = 11.89725207231041

# --- Sheet "Detailed" (sheet2): rows 39-93, column B (Price) and C (Type) ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B39").Value = 7.77743
$detailed.Range("B40").Value = 18.5855
$detailed.Range("C41").Value = "historical"
$detailed.Range("B42").Value = 64.08020999999999
$detailed.Range("C42").Value = "historical"
$detailed.Range("B43").Value = 57.31
$detailed.Range("B44").Value = 57.01318
$detailed.Range("B45").Value = 45.00272
$detailed.Range("B46").Value = 32.91642
$detailed.Range("B47").Value = 57.04922
$detailed.Range("B50").Value = 57.06
$detailed.Range("B52").Value = 57.06
$detailed.Range("B53").Value = 57.06
$detailed.Range("B54").Value = 56.97994
$detailed.Range("B55").Value = 56.97994
$detailed.Range("B56").Value = 56.98
$detailed.Range("B61").Value = 58.491
$detailed.Range("B62").Value = 64.8901
$detailed.Range("B64").Value = 26.62162
$detailed.Range("B65").Value = 35.88
$detailed.Range("B67").Value = 0.7
$detailed.Range("B69").Value = 0.7
$detailed.Range("B70").Value = 22.07
$detailed.Range("B71").Value = 24.42555
$detailed.Range("B72").Value = 22.85589
$detailed.Range("B73").Value = 26.24201
$detailed.Range("B74").Value = 28.64882
$detailed.Range("B75").Value = 36.06
$detailed.Range("B77").Value = 33.68537
$detailed.Range("B78").Value = 36.06
$detailed.Range("B80").Value = 27.46053
$detailed.Range("B81").Value = -1.17721
$detailed.Range("B82").Value = -1.96565
$detailed.Range("B83").Value = -6
$detailed.Range("B85").Value = -3.22469
$detailed.Range("B86").Value = 30.24539
$detailed.Range("B87").Value = 36.01983
$detailed.Range("B88").Value = 56.98
$detailed.Range("B90").Value = 57.09
$detailed.Range("B92").Value = 56.21361
$detailed.Range("B93").Value = 56.98
